$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6 updates: total_customers 355 -> 356, returning_customers 282 -> 283
$ws.Range("C6").Value = 356
$ws.Range("D6").Value = 283

# Derived metrics (new_customers E6 stays 73)
$ws.Range("F6").Value = 63.02895322939867
$ws.Range("G6").Value = 20.50561797752809
$ws.Range("H6").Value = 79.49438202247191
